$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-save F19 so its shared-string escaping gets normalized
# (literal "_x000D_" text picks up the correctly doubled "_x005F_x000D_" escape
# the way the rest of the sheet already does for this kind of text).
$ws.Range("F19").Value = "The website is working fire_x000D_`n"
$ws.Rows(19).AutoFit()

# Append the new feedback row
$ws.Range("D20").Value = "Lakshya Garg"
$ws.Range("E20").Value = "garglakshya635@gmail.com"
$ws.Range("F20").Value = "nigga you are doing great"
